# "Fruta / hortaliza, semanal"
# A new weekly price observation was inserted into the data table: a new
# row is inserted at row 216 (pushing the existing rows 216-272 down to
# 217-273), and the new row 216 is populated with a fresh record for the
# same market/category/variety/quality combination.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 216..272 down to 217..273, leaving a blank row 216 (style
# carried along, so the date cell keeps its date number format).
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A216").Value = 3
$ws.Range("B216").Value = "Femacal de La Calera"
$ws.Range("C216").Value = "Coquimbo"
$ws.Range("D216").Value = "2021-12-21"
$ws.Range("E216").Value = 5
$ws.Range("F216").Value = 100112031
$ws.Range("G216").Value = "Poroto verde"
$ws.Range("H216").Value = "Magnum"
$ws.Range("I216").Value = "Primera"
$ws.Range("J216").Value = 65
$ws.Range("K216").Value = 19000
$ws.Range("L216").Value = 20000
$ws.Range("M216").Value = 19538
$ws.Range("N216").Value = "`$/malla 25 kilos"
$ws.Range("O216").Value = "Provincia de Quillota"
$ws.Range("P216").Value = 782
$ws.Range("Q216").Value = 25
$ws.Range("R216").Value = "Hortaliza"
